$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 13 (B13/C13 only, no A13) held text that, after the
# edit, logically belongs elsewhere. Delete that stray row; this shifts
# every following row up by one, carrying each row's own height along.
$ws.Rows.Item(13).Delete()

# Objetivos:' answer becomes the professor's name
$ws.Range("B10").Value = '5840622 - Miguel Justino Ribeiro Barboza'
$ws.Range("C10").Value = '5840622 - Miguel Justino Ribeiro Barboza'

# 'Programa resumido:' (now row 13) becomes 'Semestral'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

# 'Programa:' (now row 15) becomes '01/01/2020'
$ws.Range("B15").Value = '01/01/2020'
$ws.Range("C15").Value = '01/01/2020'

# 'Metodo:' (now row 18) becomes the professor's name
$ws.Range("B18").Value = '5840622 - Miguel Justino Ribeiro Barboza'
$ws.Range("C18").Value = '5840622 - Miguel Justino Ribeiro Barboza'

# 'Criterio:' (now row 19) becomes the P1/P2 evaluation text
$ws.Range("B19").Value = 'Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa.'
$ws.Range("C19").Value = 'Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa.'

# 'Norma de recuperacao:' (now row 20) becomes the average formula
$ws.Range("B20").Value = 'A média do semestre será computada com base na relação:M=(P1+2P2)/3'
$ws.Range("C20").Value = 'A média do semestre será computada com base na relação:M=(P1+2P2)/3'

# 'Bibliografia:' (now row 21) becomes the recovery-grade explanation
$ws.Range("B21").Value = 'A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2'
$ws.Range("C21").Value = 'A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre.A média final, para os alunos em recuperação, será computada com base na relação abaixo:MF=(M+RC)/2'
